# Fix: Deleted five lines because of null columns.
#
# These five data rows all had DVE / VEM_bas / VEVI_bas (columns E, F, G) equal to
# zero (i.e. "null"), so they are removed entirely from the sheet:
#   row 50  -> Palmpitschroot
#   row 74  -> Tarweglutenvoer, gedroogd
#   row 151 -> Snijmais, kuil
#   row 178 -> Maisglutenvoer, vers en kuil
#   row 181 -> Tarwegistconcentraat

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete starting from the bottom-most row and working upward so that the
# row numbers of rows not yet processed remain valid.
$rowsToDelete = @(181, 178, 151, 74, 50)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
